# Add AlexaFluor 555 to the fluorophore spectra table.
# The new fluorophore is inserted as a new column at "V", pushing the
# existing Cy3..AlexaFluor750 columns (V:AI) one column to the right (W:AJ).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank column at V - shifts V:AI -> W:AJ, keeping every
# existing value/formula/style attached to its original fluorophore.
$ws.Columns("V:V").Insert()

# Populate the newly inserted column with the AlexaFluor 555 data.
$ws.Range("V7").Value = "AlexaFluor 555"   # Fluorophore name (header row)
$ws.Range("V8").Value = 0.1                # Quantum Yield
$ws.Range("V10").Value = 0.3               # Lifetime (ns)
$ws.Range("V11").Value = 155000            # M. Extinction C. (M^-1 cm^-1)

# Match the recorded selection of the edited workbook.
$ws.Range("V10").Select()
